# "Generate Report for Archive"
# The localization status report was regenerated: items that were previously
# "Ready for handoff" are now shown as "In Translation" (status text is a
# shared string, so every cell referencing it updates), and the Status
# columns shrink to fit the shorter replacement text.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: per-language status columns E (zh-cn) and F (de-de)
$wsOverview.Range("E2:F4").Value = "In Translation"

# Per-language detail sheets: Status column C
$wsZhCn.Range("C2:C4").Value = "In Translation"
$wsDeDe.Range("C2:C4").Value = "In Translation"

# Status columns narrow now that the text is shorter than "Ready for handoff"
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
